# Daily attendance processing - 2025-10-28 08:27:41
#
# Normalizes the "Recorded By" column (G) so that when the comma-separated
# list of recorders ends with the literal token "System", that token is
# moved to the front of the list (remaining recorders keep their relative
# order after it). Values that don't end with "System" (including blanks,
# or values where "System" already leads) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-NormalizedRecordedBy($val) {
    if ($null -eq $val -or $val -eq "") {
        return $val
    }

    $parts = $val.Split(",")
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $lastIdx = $parts.Length - 1
    if ($lastIdx -gt 0 -and $parts[$lastIdx] -eq "System") {
        $result = "System"
        for ($i = 0; $i -lt $lastIdx; $i++) {
            $result = $result + ", " + $parts[$i]
        }
        return $result
    }

    return $val
}

$ur = $ws.UsedRange
$firstRow = $ur.Row
$lastRow = $firstRow + $ur.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value()
    $updated = Get-NormalizedRecordedBy $current
    if ($updated -ne $current) {
        $cell.Value = $updated
    }
}
